$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number + reporting week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# Row 14 (Murder)
$ws.Range("G14").Value = 1
$ws.Range("L14").Value = -16.666666666666
$ws.Range("N14").Value = -85.294117647058

# Row 15 (Rape)
$ws.Range("C15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 26
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 85.714285714285
$ws.Range("M15").Value = 188.888888888889
$ws.Range("N15").Value = -3.703703703703

# Row 16 (Robbery)
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 36.363636363636
$ws.Range("F16").Value = 63
$ws.Range("G16").Value = 47
$ws.Range("H16").Value = 34.042553191489
$ws.Range("I16").Value = 251
$ws.Range("J16").Value = 216
$ws.Range("K16").Value = 16.203703703703
$ws.Range("L16").Value = 68.456375838926
$ws.Range("M16").Value = 29.381443298969
$ws.Range("N16").Value = -69.759036144578

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = -15.584415584415
$ws.Range("I17").Value = 406
$ws.Range("J17").Value = 375
$ws.Range("K17").Value = 8.266666666666
$ws.Range("L17").Value = 46.043165467625
$ws.Range("M17").Value = 106.091370558376
$ws.Range("N17").Value = -11.546840958605

# Row 18 (Burglary)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 135
$ws.Range("J18").Value = 155
$ws.Range("K18").Value = -12.903225806451
$ws.Range("L18").Value = 73.076923076923
$ws.Range("M18").Value = 43.617021276595
$ws.Range("N18").Value = -77.832512315270

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -22.413793103448
$ws.Range("I19").Value = 294
$ws.Range("J19").Value = 318
$ws.Range("K19").Value = -7.547169811320
$ws.Range("L19").Value = -6.369426751592
$ws.Range("M19").Value = 72.941176470588
$ws.Range("N19").Value = -4.545454545454

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 154
$ws.Range("J20").Value = 149
$ws.Range("K20").Value = 3.355704697986
$ws.Range("L20").Value = 92.5
$ws.Range("M20").Value = 214.285714285714
$ws.Range("N20").Value = -46.527777777777

# Row 21 (TOTAL)
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 216
$ws.Range("G21").Value = 237
$ws.Range("H21").Value = -8.860759493670
$ws.Range("I21").Value = 1271
$ws.Range("J21").Value = 1233
$ws.Range("K21").Value = 3.081914030819
$ws.Range("L21").Value = 38.302502720348
$ws.Range("M21").Value = 77.266387726638
$ws.Range("N21").Value = -50.254403131115

# Row 22 (Transit)
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 200
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -71.428571428571
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = -59.615384615384
$ws.Range("L22").Value = -19.230769230769
$ws.Range("M22").Value = -22.222222222222

# Row 23 (Housing)
$ws.Range("C23").Value = 15
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 25
$ws.Range("F23").Value = 52
$ws.Range("G23").Value = 46
$ws.Range("H23").Value = 13.043478260869
$ws.Range("I23").Value = 236
$ws.Range("J23").Value = 174
$ws.Range("K23").Value = 35.632183908046
$ws.Range("L23").Value = 101.709401709402
$ws.Range("M23").Value = 87.301587301587

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -16.666666666666
$ws.Range("G24").Value = 156
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 689
$ws.Range("J24").Value = 686
$ws.Range("K24").Value = 0.437317784256
$ws.Range("L24").Value = 18.793103448275
$ws.Range("M24").Value = 25.272727272727

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = -53.333333333333
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 100
$ws.Range("H25").Value = -23
$ws.Range("I25").Value = 489
$ws.Range("J25").Value = 447
$ws.Range("K25").Value = 9.395973154362
$ws.Range("L25").Value = 26.356589147286
$ws.Range("M25").Value = 1.242236024844

# Row 26 (UCR Rape*)
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = 14.285714285714
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 75
$ws.Range("L26").Value = 2.941176470588

# Row 27 (Other Sex Crimes)
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 50
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = -23.076923076923
$ws.Range("I27").Value = 56
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = 55.555555555555
$ws.Range("L27").Value = 60

# Row 28 (Shooting Vic.)
$ws.Range("C28").Value = 3
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 11
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -60.714285714285
$ws.Range("L28").Value = -60.714285714285
$ws.Range("M28").Value = -54.166666666666
$ws.Range("N28").Value = -87.058823529411

# Row 29 (Shooting Inc.)
$ws.Range("C29").Value = 3
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 50
$ws.Range("F29").Value = 4
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 11
$ws.Range("J29").Value = 22
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -54.166666666666
$ws.Range("M29").Value = -35.294117647058
$ws.Range("N29").Value = -85.333333333333

# Row 30 (Hate Crimes)
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
